# Config.xlsx update: add Spreadsheet URL / Sheet Name rows, split the
# recipient list into a proper table, and add explanatory cell comments.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 2: Email Subject (was row 2 "Subject") ---
$ws.Range("A2").Value = "Email Subject"
$ws.Range("B2").Value = "New Employe Boarding Annoucement September 2023"
$ws.Range("C2").Value = "Subject of the email to be sent."

# --- Row 3: Spreadsheet URL (new) ---
$ws.Range("A3").Value = "Spreadsheet URL"
$ws.Range("B3").Value = "https://docs.google.com/spreadsheets/d/165orIVd662-v-6BYnt8sJLENrY0tvxh3NYk63mj_Tnk/edit?usp=sharing"
$ws.Range("C3").Value = "Spreadsheet URL containing the Offering Letter Report."

# --- Row 4: Sheet Name (new) ---
$ws.Range("A4").Value = "Sheet Name"
$ws.Range("B4").Value = "Sheet1"
$ws.Range("C4").Value = "the name of the sheet in the spreadsheet."

# --- Row 5 intentionally blank, Row 6: second header for recipients table ---
$ws.Range("A6").Value = "Email Recipients"
$ws.Range("B6").Value = "Name"

# Copy the bold header look from row 1 onto the new row 6 header
$ws.Range("A1:B1").Copy()
$ws.Range("A6:B6").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

# --- Rows 7-9: individual recipients ---
$ws.Range("A7").Value = "zeerx7@gmail.com"
$ws.Range("B7").Value = "Leader"

$ws.Range("A8").Value = "pausi347@gmail.com"
$ws.Range("B8").Value = "HC TEam"

$ws.Range("A9").Value = "rosadirully5@gmail.com"
$ws.Range("B9").Value = "Leader"

# Hyperlink the email addresses (mailto:) - this also applies the
# built-in "Hyperlink" style to those cells.
$ws.Hyperlinks.Add($ws.Range("A8"), "mailto:pausi347@gmail.com")
$ws.Hyperlinks.Add($ws.Range("A9"), "mailto:rosadirully5@gmail.com")
$ws.Hyperlinks.Add($ws.Range("A7"), "mailto:zeerx7@gmail.com")

# Widen column B (no longer auto "best fit") to fit the long URL
$ws.Columns("B").ColumnWidth = 98.16666666666667

# Explanatory comments on the new recipients header cells
$ws.Range("A6").AddComment("Description:" + [char]10 + "List of destination emails to be sent." + [char]10 + "You can enter several in a vertical order")
$ws.Range("B6").AddComment("Description:" + [char]10 + "List of email owner names to be sent. ")

# Restore the originally-selected cell (shifted down with the new rows)
$ws.Range("B14").Select()
